$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet3"

# Row 1: header
$ws.Range("A1").Value = "Task ID"
$ws.Range("B1").Value = "Task Name"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Owner"
$ws.Range("E1").Value = "Effort (Days)"
$ws.Range("F1").Value = "Start Date"
$ws.Range("G1").Value = "End Date"
$ws.Range("A1:G1").Font.Bold = $true

# Row 2: section header
$ws.Range("A2").Value = "11. Test Coverage Expansion"
$ws.Range("A2").Font.Bold = $true

# Row 3: T11.1
$ws.Range("A3").Value = "T11.1"
$ws.Range("B3").Value = "Automated Test Coverage Analysis"
$ws.Range("C3").Value = "Review overall test coverage to ensure all validation aspects are automated, including validation for Silver and Gold layers."
$ws.Range("D3").Value = "Automation Engineer"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45736
$ws.Range("F3").NumberFormat = "m/d/yy"
$ws.Range("G3").Value = 45740
$ws.Range("G3").NumberFormat = "m/d/yy"

# Row 4: T11.2
$ws.Range("A4").Value = "T11.2"
$ws.Range("B4").Value = "Extend Test Coverage for Silver/Gold Layer"
$ws.Range("C4").Value = "Add test cases to the automated framework specifically for Silver and Gold layers to validate business-critical data."
$ws.Range("C4").Characters(60, 22).Font.Bold = $true
$ws.Range("D4").Value = "Automation Engineer"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45740
$ws.Range("F4").NumberFormat = "m/d/yy"
$ws.Range("G4").Value = 45743
$ws.Range("G4").NumberFormat = "m/d/yy"

# Row 5: T11.3
$ws.Range("A5").Value = "T11.3"
$ws.Range("B5").Value = "Review Test Coverage with Stakeholders"
$ws.Range("C5").Value = "Align test coverage with stakeholders to confirm all necessary validations are included and automated."
$ws.Range("D5").Value = "QA Lead"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45743
$ws.Range("F5").NumberFormat = "m/d/yy"
$ws.Range("G5").Value = 45744
$ws.Range("G5").NumberFormat = "m/d/yy"

# Row 6: T11.4
$ws.Range("A6").Value = "T11.4"
$ws.Range("B6").Value = "Add Coverage for Outbound Feeds"
$ws.Range("C6").Value = "Add automated test cases for validating outbound feeds from the ODS to ensure data integrity."
$ws.Range("C6").Characters(41, 14).Font.Bold = $true
$ws.Range("D6").Value = "Automation Engineer"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45744
$ws.Range("F6").NumberFormat = "m/d/yy"
$ws.Range("G6").Value = 45747
$ws.Range("G6").NumberFormat = "m/d/yy"

# Row 7: section header
$ws.Range("A7").Value = "12. Source DB & System Access Setup"
$ws.Range("A7").Font.Bold = $true

# Row 8: T12.1
$ws.Range("A8").Value = "T12.1"
$ws.Range("B8").Value = "Oracle DB Access"
$ws.Range("C8").Value = "Ensure access to Oracle Database for extracting legacy data to validate against ODS."
$ws.Range("C8").Characters(18, 15).Font.Bold = $true
$ws.Range("D8").Value = "Data Engineer"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45734
$ws.Range("F8").NumberFormat = "m/d/yy"
$ws.Range("G8").Value = 45735
$ws.Range("G8").NumberFormat = "m/d/yy"

# Row 9: T12.2
$ws.Range("A9").Value = "T12.2"
$ws.Range("B9").Value = "Guidewire UI Access Setup"
$ws.Range("C9").Value = "Provide access to Guidewire UI for validating data transformations and business logic."
$ws.Range("C9").Characters(19, 12).Font.Bold = $true
$ws.Range("D9").Value = "QA Tester"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45736
$ws.Range("F9").NumberFormat = "m/d/yy"
$ws.Range("G9").Value = 45737
$ws.Range("G9").NumberFormat = "m/d/yy"

# Row 10: T12.3
$ws.Range("A10").Value = "T12.3"
$ws.Range("B10").Value = "S3 Access for Staging"
$ws.Range("C10").Value = "Ensure access to S3 buckets for data validation between raw and ODS layers."
$ws.Range("C10").Characters(18, 10).Font.Bold = $true
$ws.Range("D10").Value = "Data Engineer"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45737
$ws.Range("F10").NumberFormat = "m/d/yy"
$ws.Range("G10").Value = 45738
$ws.Range("G10").NumberFormat = "m/d/yy"

# Row 11: T12.4
$ws.Range("A11").Value = "T12.4"
$ws.Range("B11").Value = "Jump Server Access"
$ws.Range("C11").Value = "Ensure access to Jump Server for secure connections during SIT testing."
$ws.Range("C11").Characters(18, 11).Font.Bold = $true
$ws.Range("D11").Value = "Data Engineer"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45738
$ws.Range("F11").NumberFormat = "m/d/yy"
$ws.Range("G11").Value = 45739
$ws.Range("G11").NumberFormat = "m/d/yy"

# Row 12: section header
$ws.Range("A12").Value = "13. DevOps Setup & CI/CD Integration"
$ws.Range("A12").Font.Bold = $true

# Row 13: T13.1
$ws.Range("A13").Value = "T13.1"
$ws.Range("B13").Value = "DevOps Feasibility for CI/CD"
$ws.Range("C13").Value = "Investigate CI/CD integration for automated testing to ensure continuous deployment pipelines are effective."
$ws.Range("C13").Characters(13, 17).Font.Bold = $true
$ws.Range("D13").Value = "DevOps Engineer"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 45738
$ws.Range("F13").NumberFormat = "m/d/yy"
$ws.Range("G13").Value = 45741
$ws.Range("G13").NumberFormat = "m/d/yy"

# Row 14: T13.2
$ws.Range("A14").Value = "T13.2"
$ws.Range("B14").Value = "Repository Creation in DevOps"
$ws.Range("C14").Value = "Set up a GitHub/DevOps repository to store all automation scripts, configurations, and test artifacts."
$ws.Range("C14").Characters(10, 24).Font.Bold = $true
$ws.Range("D14").Value = "DevOps Engineer"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 45741
$ws.Range("F14").NumberFormat = "m/d/yy"
$ws.Range("G14").Value = 45742
$ws.Range("G14").NumberFormat = "m/d/yy"

# Row 15: T13.3
$ws.Range("A15").Value = "T13.3"
$ws.Range("B15").Value = "Set up CI/CD Pipelines"
$ws.Range("C15").Value = "Configure CI/CD pipelines in Azure DevOps for automated testing integration with deployment processes."
$ws.Range("C15").Characters(11, 15).Font.Bold = $true
$ws.Range("C15").Characters(30, 12).Font.Bold = $true
$ws.Range("D15").Value = "DevOps Engineer"
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 45742
$ws.Range("F15").NumberFormat = "m/d/yy"
$ws.Range("G15").Value = 45746
$ws.Range("G15").NumberFormat = "m/d/yy"

# Row 16: section header
$ws.Range("A16").Value = "14. Data Integration Testing"
$ws.Range("A16").Font.Bold = $true

# Row 17: T14.1
$ws.Range("A17").Value = "T14.1"
$ws.Range("B17").Value = "Data Integration Testing (S3 ↔ ODS)"
$ws.Range("C17").Value = "Ensure automated testing of data from S3 staging to ODS and validate data integrity and transformations."
$ws.Range("C17").Characters(39, 10).Font.Bold = $true
$ws.Range("C17").Characters(53, 3).Font.Bold = $true
$ws.Range("D17").Value = "QA Tester"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 45746
$ws.Range("F17").NumberFormat = "m/d/yy"
$ws.Range("G17").Value = 45749
$ws.Range("G17").NumberFormat = "m/d/yy"

# Row 18: T14.2
$ws.Range("A18").Value = "T14.2"
$ws.Range("B18").Value = "Data Integration Testing (PostgreSQL ↔ ODS)"
$ws.Range("C18").Value = "Validate data from Read Replica PostgreSQL to ODS for correctness and consistency."
$ws.Range("C18").Characters(20, 23).Font.Bold = $true
$ws.Range("C18").Characters(47, 3).Font.Bold = $true
$ws.Range("D18").Value = "QA Tester"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 45749
$ws.Range("F18").NumberFormat = "m/d/yy"
$ws.Range("G18").Value = 45752
$ws.Range("G18").NumberFormat = "m/d/yy"

# Row 19: T14.3
$ws.Range("A19").Value = "T14.3"
$ws.Range("B19").Value = "Data Integration Testing (S3 ↔ ODS_VIEWS)"
$ws.Range("C19").Value = "Ensure consistency and correctness of data between S3 staging and ODS_VIEWS tables."
$ws.Range("C19").Characters(52, 10).Font.Bold = $true
$ws.Range("C19").Characters(67, 9).Font.Bold = $true
$ws.Range("D19").Value = "QA Tester"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 45752
$ws.Range("F19").NumberFormat = "m/d/yy"
$ws.Range("G19").Value = 45755
$ws.Range("G19").NumberFormat = "m/d/yy"

# Row 20: section header
$ws.Range("A20").Value = "15. Automation Maintenance & Monitoring"
$ws.Range("A20").Font.Bold = $true

# Row 21: T15.1
$ws.Range("A21").Value = "T15.1"
$ws.Range("B21").Value = "Automation Script Maintenance"
$ws.Range("C21").Value = "Set up a maintenance plan for updating scripts based on any changes in source data or ODS schemas."
$ws.Range("C21").Characters(10, 16).Font.Bold = $true
$ws.Range("D21").Value = "Automation Engineer"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 45755
$ws.Range("F21").NumberFormat = "m/d/yy"
$ws.Range("G21").Value = 45758
$ws.Range("G21").NumberFormat = "m/d/yy"

# Row 22: T15.2
$ws.Range("A22").Value = "T15.2"
$ws.Range("B22").Value = "Monitor Automation Results"
$ws.Range("C22").Value = "Regularly monitor and ensure that automation results are generated and reports are shared with stakeholders."
$ws.Range("C22").Characters(35, 18).Font.Bold = $true
$ws.Range("D22").Value = "QA Tester"
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 45758
$ws.Range("F22").NumberFormat = "m/d/yy"
$ws.Range("G22").Value = 45760
$ws.Range("G22").NumberFormat = "m/d/yy"

# Row 23: section header
$ws.Range("A23").Value = "16. Final Testing & SIT Execution"
$ws.Range("A23").Font.Bold = $true

# Row 24: T16.1
$ws.Range("A24").Value = "T16.1"
$ws.Range("B24").Value = "SIT Execution & Test Suite Running"
$ws.Range("C24").Value = "Execute full SIT testing suite, including automated and manual validation, for end-to-end integration."
$ws.Range("C24").Characters(9, 16).Font.Bold = $true
$ws.Range("D24").Value = "QA Tester"
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 45788
$ws.Range("F24").NumberFormat = "m/d/yy"
$ws.Range("G24").Value = 45793
$ws.Range("G24").NumberFormat = "m/d/yy"

# Row 25: T16.2
$ws.Range("A25").Value = "T16.2"
$ws.Range("B25").Value = "Defect Tracking & Re-testing"
$ws.Range("C25").Value = "Ensure that all defects are logged, tracked in Jira, and re-tested once fixed."
$ws.Range("C25").Characters(13, 11).Font.Bold = $true
$ws.Range("D25").Value = "QA Tester / Automation Engineer"
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 45793
$ws.Range("F25").NumberFormat = "m/d/yy"
$ws.Range("G25").Value = 45797
$ws.Range("G25").NumberFormat = "m/d/yy"

# Row 26: T16.3
$ws.Range("A26").Value = "T16.3"
$ws.Range("B26").Value = "Final SIT Sign-off"
$ws.Range("C26").Value = "Final review, sign-off from stakeholders, and confirmation of test completion."
$ws.Range("D26").Value = "Project Lead"
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 45797
$ws.Range("F26").NumberFormat = "m/d/yy"
$ws.Range("G26").Value = 45799
$ws.Range("G26").NumberFormat = "m/d/yy"

# Selection to match target sheetView
[void]$ws.Range("A1:G26").Select()
